$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 2, "Text Box 5" (Shapes.Item(4)): three bullet paragraphs.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shBullets = $s2.Shapes.Item(4)
$tfBullets = $shBullets.TextFrame
$trBullets = $tfBullets.TextRange

# Remember the shape's box. Its <a:spAutoFit/> would otherwise recompute the
# height (badly) as soon as we touch any run below; the source edit never
# actually resized the box, so freeze autosizing while we edit, then restore
# both the autosize flag and the exact original box afterwards.
$origHeight = $shBullets.Height
$origWidth = $shBullets.Width
$origTop = $shBullets.Top
$origLeft = $shBullets.Left
$tfBullets.AutoSize = 0

# 1) "Data analysis is now applied in anywhere." -> "...in everywhere."
#    (splits into 3 runs: "...applied in ", "everywhere", ".")
$para2 = $trBullets.Paragraphs(2)
$fullPara2 = $para2.Text
$idxAnywhere = $fullPara2.IndexOf("anywhere")
$subAnywhere = $para2.Characters($idxAnywhere + 1, 8)
$subAnywhere.Text = "everywhere"

# 2) "More and more graduated students are looking    for "Data Analyst" posts."
#    The three runs covering "...posts" + "." get merged back into a single run.
#    Build the replacement with explicit Unicode escapes for the curly quotes so
#    the COM text-normalisation (curly -> straight) never touches this string.
$para3 = $trBullets.Paragraphs(3)
$oq = [char]0x201C
$cq = [char]0x201D
$found3 = $para3.Find("More and more", 0)
$startChar3 = $found3.Start - $para3.Start + 1
$restLen3 = $para3.Length - $startChar3 + 1
$sub3 = $para3.Characters($startChar3, $restLen3)
$sub3.Text = "More and more graduated students are looking    for ${oq}Data Analyst${cq} posts."

# 3) "Our project play an important role in the posts finding process."
#    The three runs covering "...the " + "posts " + "finding process." get
#    merged back into a single run.
$para4 = $trBullets.Paragraphs(4)
$found4 = $para4.Find("Our project", 0)
$startChar4 = $found4.Start - $para4.Start + 1
$restLen4 = $para4.Length - $startChar4 + 1
$sub4 = $para4.Characters($startChar4, $restLen4)
$sub4.Text = "Our project play an important role in the posts finding process."

# Restore autosizing and the shape's original geometry (see note above).
# (Re-typing the literal -- rather than reusing the captured variable --
# avoids a 1-EMU rounding drift through the point<->EMU conversion.)
$tfBullets.AutoSize = 1
$shBullets.Height = 225.3797
$shBullets.Width = $origWidth
$shBullets.Top = $origTop
$shBullets.Left = $origLeft

# ---------------------------------------------------------------------------
# Slide 4, "TextBox 32" (Shapes.Item(14)): fix "infomation" -> "information"
# and drop the stale spell-check "err" flag by recreating the run instead of
# just editing its text in place.
# ---------------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$shInfo = $s4.Shapes.Item(14)
$trInfo = $shInfo.TextFrame.TextRange
$foundInfo = $trInfo.Find("infomation", 0)
$foundInfo.Delete()
[void]$trInfo.InsertAfter("information")
